$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(43, 8).Value = 2627.2222
$ws.Cells.Item(43, 10).Value = 1997.5
$ws.Cells.Item(43, 12).Value = 1997.5
$ws.Cells.Item(43, 14).Value = -2135.5
$ws.Cells.Item(126, 8).Value = 46766
$ws.Cells.Item(126, 10).Value = 46766
$ws.Cells.Item(126, 12).Value = 46766
$ws.Cells.Item(126, 14).Value = -56646
$ws.Cells.Item(129, 8).Value = 1458.3658
$ws.Cells.Item(129, 9).Value = 1377.8
$ws.Cells.Item(129, 10).Value = 1484.3549
$ws.Cells.Item(129, 11).Value = 4133.4
$ws.Cells.Item(129, 12).Value = 4453.0647
$ws.Cells.Item(129, 13).Value = 866.6000000000004
$ws.Cells.Item(129, 14).Value = -14453.0647
$ws.Cells.Item(130, 8).Value = 49772
$ws.Cells.Item(130, 10).Value = 49772
$ws.Cells.Item(130, 12).Value = 49772
$ws.Cells.Item(130, 14).Value = -59812
$ws.Cells.Item(133, 8).Value = 53816.332
$ws.Cells.Item(133, 10).Value = 53816.332
$ws.Cells.Item(133, 12).Value = 53816.332
$ws.Cells.Item(133, 14).Value = -63936.332
$ws.Cells.Item(138, 8).Value = 1402.41
$ws.Cells.Item(138, 9).Value = 690.44684
$ws.Cells.Item(138, 10).Value = 2033.7736
$ws.Cells.Item(138, 11).Value = 2071.34052
$ws.Cells.Item(138, 12).Value = 6101.3208
$ws.Cells.Item(138, 13).Value = 3068.65948
$ws.Cells.Item(138, 14).Value = -16381.3208

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 1832.6666
$ws.Cells.Item(2, 9).Value = 2002.2
$ws.Cells.Item(2, 10).Value = 985
$ws.Cells.Item(2, 11).Value = 2002.2
$ws.Cells.Item(2, 12).Value = 985
$ws.Cells.Item(2, 13).Value = -1889.2
$ws.Cells.Item(2, 14).Value = -1211
$ws.Cells.Item(61, 8).Value = 1377.1562
$ws.Cells.Item(61, 9).Value = 1002.4
$ws.Cells.Item(61, 10).Value = 2715.5715
$ws.Cells.Item(61, 11).Value = 1002.4
$ws.Cells.Item(61, 12).Value = 2715.5715
$ws.Cells.Item(61, 13).Value = -790.4
$ws.Cells.Item(61, 14).Value = -3139.5715
$ws.Cells.Item(116, 8).Value = 1832.6666
$ws.Cells.Item(116, 9).Value = 2002.2
$ws.Cells.Item(116, 10).Value = 985
$ws.Cells.Item(116, 11).Value = 2002.2
$ws.Cells.Item(116, 12).Value = 985
$ws.Cells.Item(116, 13).Value = 291.8
$ws.Cells.Item(116, 14).Value = -5573
$ws.Cells.Item(122, 8).Value = 1731.5
$ws.Cells.Item(122, 9).Value = 1602.25
$ws.Cells.Item(122, 10).Value = 1990
$ws.Cells.Item(122, 11).Value = 4806.75
$ws.Cells.Item(122, 12).Value = 5970
$ws.Cells.Item(122, 13).Value = -2356.75
$ws.Cells.Item(122, 14).Value = -10870
$ws.Cells.Item(125, 8).Value = 45283.832
$ws.Cells.Item(125, 10).Value = 45283.832
$ws.Cells.Item(125, 12).Value = 45283.832
$ws.Cells.Item(125, 14).Value = -55123.832
$ws.Cells.Item(131, 8).Value = 47037.332
$ws.Cells.Item(131, 10).Value = 47037.332
$ws.Cells.Item(131, 12).Value = 47037.332
$ws.Cells.Item(131, 14).Value = -57117.332
$ws.Cells.Item(136, 8).Value = 1377.1562
$ws.Cells.Item(136, 9).Value = 1002.4
$ws.Cells.Item(136, 10).Value = 2715.5715
$ws.Cells.Item(136, 11).Value = 3007.2
$ws.Cells.Item(136, 12).Value = 8146.7145
$ws.Cells.Item(136, 13).Value = -457.1999999999998
$ws.Cells.Item(136, 14).Value = -13246.7145
$ws.Cells.Item(138, 8).Value = 53400
$ws.Cells.Item(138, 10).Value = 53400
$ws.Cells.Item(138, 12).Value = 53400
$ws.Cells.Item(138, 14).Value = -63680

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 1832.6666
$ws.Cells.Item(3, 9).Value = 2002.2
$ws.Cells.Item(3, 10).Value = 985
$ws.Cells.Item(3, 11).Value = 2002.2
$ws.Cells.Item(3, 12).Value = 985
$ws.Cells.Item(3, 13).Value = -1888.2
$ws.Cells.Item(3, 14).Value = -1213
$ws.Cells.Item(7, 8).Value = 650
$ws.Cells.Item(7, 9).Value = 650
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 650
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = -537
$ws.Cells.Item(7, 14).ClearContents()
$ws.Cells.Item(94, 8).Value = 1920
$ws.Cells.Item(94, 9).Value = 1920
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 1920
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = -1469
$ws.Cells.Item(94, 14).ClearContents()
$ws.Cells.Item(105, 8).Value = 3756.7778
$ws.Cells.Item(105, 9).Value = 3160
$ws.Cells.Item(105, 10).Value = 4502.75
$ws.Cells.Item(105, 11).Value = 3160
$ws.Cells.Item(105, 12).Value = 4502.75
$ws.Cells.Item(105, 13).Value = -1413
$ws.Cells.Item(105, 14).Value = -7996.75
$ws.Cells.Item(122, 8).Value = 40577.6
$ws.Cells.Item(122, 10).Value = 40577.6
$ws.Cells.Item(122, 12).Value = 40577.6
$ws.Cells.Item(122, 14).Value = -50377.6
$ws.Cells.Item(124, 8).Value = 50992
$ws.Cells.Item(124, 10).Value = 50992
$ws.Cells.Item(124, 12).Value = 50992
$ws.Cells.Item(124, 14).Value = -60812

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 4522.9595
$ws.Cells.Item(31, 9).Value = 2030.1177
$ws.Cells.Item(31, 10).Value = 5039.768
$ws.Cells.Item(31, 11).Value = 2030.1177
$ws.Cells.Item(31, 12).Value = 5039.768
$ws.Cells.Item(31, 13).Value = -1735.1177
$ws.Cells.Item(31, 14).Value = -5629.768
$ws.Cells.Item(34, 8).Value = 4522.9595
$ws.Cells.Item(34, 9).Value = 2030.1177
$ws.Cells.Item(34, 10).Value = 5039.768
$ws.Cells.Item(34, 11).Value = 2030.1177
$ws.Cells.Item(34, 12).Value = 5039.768
$ws.Cells.Item(34, 13).Value = -1828.1177
$ws.Cells.Item(34, 14).Value = -5443.768
$ws.Cells.Item(100, 8).Value = 33330.4
$ws.Cells.Item(100, 10).Value = 33330.4
$ws.Cells.Item(100, 12).Value = 33330.4
$ws.Cells.Item(100, 14).Value = -35494.4
$ws.Cells.Item(134, 8).Value = 23360.572
$ws.Cells.Item(134, 9).Value = 1036.4546
$ws.Cells.Item(134, 10).Value = 176838.88
$ws.Cells.Item(134, 11).Value = 3109.3638
$ws.Cells.Item(134, 12).Value = 530516.64
$ws.Cells.Item(134, 13).Value = -574.3638000000001
$ws.Cells.Item(134, 14).Value = -535586.64
$ws.Cells.Item(137, 8).Value = 46499.92
$ws.Cells.Item(137, 10).Value = 46499.92
$ws.Cells.Item(137, 12).Value = 46499.92
$ws.Cells.Item(137, 14).Value = -56699.92

# Sheet 5
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(4, 8).Value = 1598
$ws.Cells.Item(4, 9).Value = 95.40000000000001
$ws.Cells.Item(4, 10).Value = 2850.1667
$ws.Cells.Item(4, 11).Value = 286.2
$ws.Cells.Item(4, 12).Value = 8550.500100000001
$ws.Cells.Item(4, 13).Value = -174.2
$ws.Cells.Item(4, 14).Value = -8774.500100000001
$ws.Cells.Item(5, 8).Value = 4247.7188
$ws.Cells.Item(5, 9).Value = 9735.637000000001
$ws.Cells.Item(5, 10).Value = 1373.0952
$ws.Cells.Item(5, 11).Value = 29206.911
$ws.Cells.Item(5, 12).Value = 4119.2856
$ws.Cells.Item(5, 13).Value = -29094.911
$ws.Cells.Item(5, 14).Value = -4343.2856
$ws.Cells.Item(39, 8).Value = 5125
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 5125
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 12).Value = 15375
$ws.Cells.Item(39, 13).ClearContents()
$ws.Cells.Item(39, 14).Value = -15963
$ws.Cells.Item(132, 8).Value = 2024.1177
$ws.Cells.Item(132, 9).Value = 1345.1875
$ws.Cells.Item(132, 11).Value = 12106.6875
$ws.Cells.Item(132, 13).Value = -9576.6875
$ws.Cells.Item(135, 8).Value = 4247.7188
$ws.Cells.Item(135, 9).Value = 9735.637000000001
$ws.Cells.Item(135, 10).Value = 1373.0952
$ws.Cells.Item(135, 11).Value = 87620.73300000001
$ws.Cells.Item(135, 12).Value = 12357.8568
$ws.Cells.Item(135, 13).Value = -85085.73300000001
$ws.Cells.Item(135, 14).Value = -17427.8568
$ws.Cells.Item(140, 8).Value = 162021.23
$ws.Cells.Item(140, 9).Value = 201479.95
$ws.Cells.Item(140, 10).Value = 4186.4
$ws.Cells.Item(140, 11).Value = 604439.8500000001
$ws.Cells.Item(140, 12).Value = 12559.2
$ws.Cells.Item(140, 13).Value = -599259.8500000001
$ws.Cells.Item(140, 14).Value = -22919.2

# Sheet 6
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(113, 8).Value = 1163.1177
$ws.Cells.Item(113, 9).Value = 1069.091
$ws.Cells.Item(113, 10).Value = 1335.5
$ws.Cells.Item(113, 11).Value = 1069.091
$ws.Cells.Item(113, 12).Value = 1335.5
$ws.Cells.Item(113, 13).Value = 1100.909
$ws.Cells.Item(113, 14).Value = -5675.5
$ws.Cells.Item(122, 8).Value = 1087.5
$ws.Cells.Item(122, 9).Value = 1160
$ws.Cells.Item(122, 10).Value = 966.6667
$ws.Cells.Item(122, 11).Value = 3480
$ws.Cells.Item(122, 12).Value = 2900.0001
$ws.Cells.Item(122, 13).Value = -1030
$ws.Cells.Item(122, 14).Value = -7800.0001
$ws.Cells.Item(124, 8).Value = 41511.668
$ws.Cells.Item(124, 10).Value = 41511.668
$ws.Cells.Item(124, 12).Value = 41511.668
$ws.Cells.Item(124, 14).Value = -51331.668
$ws.Cells.Item(126, 8).Value = 6427.231
$ws.Cells.Item(126, 9).Value = 7668.8423
$ws.Cells.Item(126, 10).Value = 3057.1428
$ws.Cells.Item(126, 11).Value = 23006.5269
$ws.Cells.Item(126, 12).Value = 9171.428400000001
$ws.Cells.Item(126, 13).Value = -20536.5269
$ws.Cells.Item(126, 14).Value = -14111.4284
$ws.Cells.Item(135, 8).Value = 46850
$ws.Cells.Item(135, 10).Value = 46850
$ws.Cells.Item(135, 12).Value = 46850
$ws.Cells.Item(135, 14).Value = -56990

# Sheet 7
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 2295.8857
$ws.Cells.Item(7, 9).Value = 1794.6786
$ws.Cells.Item(7, 10).Value = 4300.7144
$ws.Cells.Item(7, 11).Value = 1794.6786
$ws.Cells.Item(7, 12).Value = 4300.7144
$ws.Cells.Item(7, 13).Value = -1682.6786
$ws.Cells.Item(7, 14).Value = -4524.7144
$ws.Cells.Item(98, 8).Value = 37325.332
$ws.Cells.Item(98, 10).Value = 37325.332
$ws.Cells.Item(98, 12).Value = 37325.332
$ws.Cells.Item(98, 14).Value = -43315.332
$ws.Cells.Item(126, 8).Value = 2295.8857
$ws.Cells.Item(126, 9).Value = 1794.6786
$ws.Cells.Item(126, 10).Value = 4300.7144
$ws.Cells.Item(126, 11).Value = 5384.0358
$ws.Cells.Item(126, 12).Value = 12902.1432
$ws.Cells.Item(126, 13).Value = -2914.0358
$ws.Cells.Item(126, 14).Value = -17842.1432
$ws.Cells.Item(137, 8).Value = 41316.668
$ws.Cells.Item(137, 10).Value = 41316.668
$ws.Cells.Item(137, 12).Value = 41316.668
$ws.Cells.Item(137, 14).Value = -51516.668

# Sheet 8
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(4, 8).Value = 3001000
$ws.Cells.Item(4, 9).Value = 6000000
$ws.Cells.Item(4, 10).Value = 2000
$ws.Cells.Item(4, 11).Value = 6000000
$ws.Cells.Item(4, 12).Value = 2000
$ws.Cells.Item(4, 13).Value = -5999887
$ws.Cells.Item(4, 14).Value = -2226
$ws.Cells.Item(103, 8).Value = 49594
$ws.Cells.Item(103, 10).Value = 49594
$ws.Cells.Item(103, 12).Value = 49594
$ws.Cells.Item(103, 14).Value = -51938
$ws.Cells.Item(122, 8).Value = 1361197.8
$ws.Cells.Item(122, 9).Value = 1786284.6
$ws.Cells.Item(122, 10).Value = 920
$ws.Cells.Item(122, 11).Value = 5358853.800000001
$ws.Cells.Item(122, 12).Value = 2760
$ws.Cells.Item(122, 13).Value = -5356403.800000001
$ws.Cells.Item(122, 14).Value = -7660
$ws.Cells.Item(126, 8).Value = 806.5
$ws.Cells.Item(126, 9).Value = 815.9167
$ws.Cells.Item(126, 10).Value = 750
$ws.Cells.Item(126, 11).Value = 2447.7501
$ws.Cells.Item(126, 12).Value = 2250
$ws.Cells.Item(126, 13).Value = 22.2498999999998
$ws.Cells.Item(126, 14).Value = -7190
$ws.Cells.Item(131, 8).Value = 49215.668
$ws.Cells.Item(131, 10).Value = 49215.668
$ws.Cells.Item(131, 12).Value = 49215.668
$ws.Cells.Item(131, 14).Value = -59295.668
$ws.Cells.Item(139, 8).Value = 57800
$ws.Cells.Item(139, 10).Value = 57800
$ws.Cells.Item(139, 12).Value = 57800
$ws.Cells.Item(139, 14).Value = -68080
